$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting the rest of the table down
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row 2 with the new customer record
$ws.Cells.Item(2, 1).Value = "KH"
$ws.Cells.Item(2, 2).Value = 395
$ws.Cells.Item(2, 3).Value = "ngọc hân"
$ws.Cells.Item(2, 4).Value = "SÓC TRĂNG"
$ws.Cells.Item(2, 5).Value = $null
$ws.Cells.Item(2, 6).Value = $null
$ws.Cells.Item(2, 7).Value = $null
$ws.Cells.Item(2, 8).Value = $null
$ws.Cells.Item(2, 9).Value = 35000000
$ws.Cells.Item(2, 10).Value = 8000000
